$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

$ws.Range('D2').Value = '29.398.28'
$ws.Range('E2').Value = '  -0.12%  '

$ws.Range('D3').Value = '1.844.74'
$ws.Range('E3').Value = '  -0.32%  '

$ws.Range('E4').Value = '  +0.03%  '

Set-TextValue $ws 'D5' '239.04'
$ws.Range('E5').Value = '  -0.83%  '

Set-TextValue $ws 'D6' '0.6316'
$ws.Range('E6').Value = '  -0.33%  '

$ws.Range('E7').Value = '  +0.04%  '

Set-TextValue $ws 'D8' '0.07538'
$ws.Range('E8').Value = '  -0.43%  '

Set-TextValue $ws 'D9' '0.2931'
$ws.Range('E9').Value = '  -1.39%  '

Set-TextValue $ws 'D10' '24.55'
$ws.Range('E10').Value = '  -0.54%  '

Set-TextValue $ws 'D11' '0.07713'
$ws.Range('E11').Value = '  -0.16%  '

$ws.Range('D12').Value = '1.841.05'
$ws.Range('E12').Value = '  -7.25%  '

Set-TextValue $ws 'D13' '5.002'
$ws.Range('E13').Value = '  +0.04%  '

Set-TextValue $ws 'D14' '0.6800'

$ws.Range('E15').Value = '  +4.87%  '

Set-TextValue $ws 'D16' '83.33'
$ws.Range('E16').Value = '  +0.28%  '

$ws.Range('D17').Value = '2.087.56'
$ws.Range('E17').Value = '  -7.81%  '

Set-TextValue $ws 'D18' '6.173'
$ws.Range('E18').Value = '  -0.73%  '

$ws.Range('D19').Value = '29.430.26'
$ws.Range('E19').Value = '  -0.14%  '

Set-TextValue $ws 'D20' '228.77'

$ws.Range('E21').Value = '  -0.67%  '

$ws.Range('E23').Value = '  -1.99%  '

$ws.Range('E24').Value = '  +0.06%  '

Set-TextValue $ws 'D25' '156.76'
$ws.Range('E25').Value = '  +0.57%  '

Set-TextValue $ws 'D26' '0.1394'
$ws.Range('E26').Value = '  +0.33%  '

Set-TextValue $ws 'D27' '8.357'
$ws.Range('E27').Value = '  -0.78%  '

Set-TextValue $ws 'D28' '17.60'
$ws.Range('E28').Value = '  -0.61%  '

$ws.Range('E29').Value = '  -0.96%  '

Set-TextValue $ws 'D30' '1.288'
$ws.Range('E30').Value = '  +2.19%  '

Set-TextValue $ws 'D31' '0.05629'
$ws.Range('E31').Value = '  -3.41%  '

Set-TextValue $ws 'D32' '4.106'
$ws.Range('E32').Value = '  -0.85%  '

Set-TextValue $ws 'D33' '4.026'
$ws.Range('E33').Value = '  -0.02%  '

Set-TextValue $ws 'D34' '1.848'
$ws.Range('E34').Value = '  -0.74%  '

Set-TextValue $ws 'D35' '1.157'
$ws.Range('E35').Value = '  -0.20%  '

Set-TextValue $ws 'D36' '0.7130'
$ws.Range('E36').Value = '  -0.65%  '

Set-TextValue $ws 'D37' '2.591'
$ws.Range('E37').Value = '  -0.10%  '

$ws.Range('D38').Value = '1.247.67'
$ws.Range('E38').Value = '  -0.47%  '

Set-TextValue $ws 'D39' '0.01810'
$ws.Range('E39').Value = '  +0.11%  '

Set-TextValue $ws 'D40' '2.770'
$ws.Range('E40').Value = '  -1.01%  '

Set-TextValue $ws 'D41' '6.374'
$ws.Range('E41').Value = '  +4.61%  '

Set-TextValue $ws 'D42' '0.9022'
$ws.Range('E42').Value = '  -0.12%  '

$ws.Range('E43').Value = '  +0.04%  '

Set-TextValue $ws 'D44' '101.68'
$ws.Range('E44').Value = '  -0.05%  '

Set-TextValue $ws 'D45' '65.85'

$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws 'D46' '0.00000000119'
$ws.Range('E46').Value = '  +1.01%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws 'D47' '7.107'
$ws.Range('E47').Value = '  -1.74%  '

$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws 'D48' '0.3997'
$ws.Range('E48').Value = '  -0.68%  '

Set-TextValue $ws 'D49' '1.673'
$ws.Range('E49').Value = '  -1.20%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D50' '8.924'
$ws.Range('E50').Value = '  -2.89%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 'D51' '0.1123'
$ws.Range('E51').Value = '  -0.43%  '
